# Atualizei dados bibi e add
# - Troca os rotulos das lojas "Bibi Cell Ponta Negra" (linha 4) e
#   "Bibi Cell Vieiralves" (linha 5), que fica refletido tambem na troca
#   dos valores diarios (B:G) entre essas duas linhas.
# - Adiciona faturamento do dia 7 (coluna H) para todas as lojas.
# - Atualiza os totais (coluna AG) de cada linha e a linha "total" (linha 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Linha 4: passa a ser "Bibi Cell Vieiralves" ---
$ws.Range("A4").Value = "Bibi Cell Vieiralves"
$ws.Range("B4").Value = 3638
$ws.Range("C4").Value = 3280.25
$ws.Range("D4").Value = 5521.8
$ws.Range("E4").Value = 2850
$ws.Range("F4").Value = 4180
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 4097.5
$ws.Range("AG4").Value = 23567.55

# --- Linha 5: passa a ser "Bibi Cell Ponta Negra" ---
$ws.Range("A5").Value = "Bibi Cell Ponta Negra"
$ws.Range("B5").Value = 4535.01
$ws.Range("C5").Value = 2416
$ws.Range("D5").Value = 2403.81
$ws.Range("E5").Value = 5469
$ws.Range("F5").Value = 3040.95
$ws.Range("G5").Value = 1795.75
$ws.Range("H5").Value = 2017.01
$ws.Range("AG5").Value = 21677.53

# --- Linha 2: Bibi Cell Mundi - novo valor do dia 7 e total ---
$ws.Range("H2").Value = 10126.16
$ws.Range("AG2").Value = 57060.47

# --- Linha 3: Bibi Cell Manauara - novo valor do dia 7 e total ---
$ws.Range("H3").Value = 3677
$ws.Range("AG3").Value = 25540.9

# --- Linha 6: total geral - novo valor do dia 7 e total ---
$ws.Range("H6").Value = 19917.67
$ws.Range("AG6").Value = 127846.45
